$d = $word.ActiveDocument

# 1. Remove the old _GoBack bookmark (it previously sat after "MP73010")
$d.Bookmarks("_GoBack").Delete()

# 2. Merge the ">>  >  your" / " stuff after this line >>>" runs (with proofErr
#    marks) into a single run reading ">>>  your stuff after this line >>>"
$findRepl = $d.Content
$findRepl.Find.Execute(">>>  your stuff after this line >>>", $false, $false, $false, $false, $false, $true, 1, $false, ">>>  your stuff after this line >>>", 2)

# 3. Insert a brand-new paragraph after "Ben changing things up!" with the
#    student's name/ID text, and re-create the _GoBack bookmark at the very
#    end of that new paragraph.
$benPara = $d.Paragraphs(5)
$afterBen = $benPara.Range
$afterBen.Collapse(0)
$afterBen.InsertParagraphAfter()

$newPara = $d.Paragraphs(6)
$newPara.Range.Text = "Zhecheng Cao, Student ID: 22863311, Hello World~"

$markerRange = $d.Content
$markerRange.Find.Execute("~", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $markerRange)
$markerRange.Text = ""
